# Add season-record columns (Wins / Losses / Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (style index 1: bold, centered, bordered)
# from the existing "Salary" header cell (AA1) onto the three new header
# cells, then set their text.
$ws.Range("AA1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Every player row (2-45) gets the same team season record.
$wins = 78
$losses = 84
$ties = 0

for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 29).Value = $wins    # column AC
    $ws.Cells.Item($row, 30).Value = $losses  # column AD
    $ws.Cells.Item($row, 31).Value = $ties    # column AE
}

Write-Host "Added Wins/Losses/Ties columns (AC:AE) for rows 1-45"
